$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - TSM
$ws.Range("D2").Value = 292.29
$ws.Range("E2").Value = 59.6
$ws.Range("F2").Value = 0.8
$ws.Range("K2").Value = 59.3
$ws.Range("N2").Value = 53.62998959737769

# Row 3 - ASML
$ws.Range("D3").Value = 1110.81
$ws.Range("E3").Value = 62.8
$ws.Range("F3").Value = 6.71
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 59.3
$ws.Range("N3").Value = 53.62998959737769

# Row 4 - QCOM
$ws.Range("D4").Value = 174.19
$ws.Range("E4").Value = 49.6
$ws.Range("F4").Value = 5.48
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 40
$ws.Range("K4").Value = 50.1
$ws.Range("N4").Value = 53.62998959737769

# Row 5 - NVDA
$ws.Range("D5").Value = 183.79
$ws.Range("E5").Value = 46.9
$ws.Range("F5").Value = 1.96
$ws.Range("K5").Value = 49.1
$ws.Range("N5").Value = 53.62998959737769

# Row 6 - AMD
$ws.Range("D6").Value = 218.1
$ws.Range("E6").Value = 32.5
$ws.Range("F6").Value = 1.8
$ws.Range("H6").Value = 56
$ws.Range("K6").Value = 46.1
$ws.Range("N6").Value = 53.62998959737769
